# Rename header columns to short machine-friendly names.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case every word in the Spanish place names held in columns A and B
# (rows 2..1705), e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga".
# Build the new word with [string]::Concat (NOT "+") because this host
# coerces "+" on two all-numeric strings (e.g. "0" and "8") into integer
# addition, which would corrupt tokens like "08" inside place names such
# as "San Juan Mixtepec - Distr. 08 -".
function TitleCaseWords($s) {
    $words = $s.Split(" ")
    $out = @()
    foreach ($w in $words) {
        if ($w.Length -gt 0) {
            $first = $w.Substring(0, 1).ToUpper()
            $rest = $w.Substring(1)
            $out += [string]::Concat($first, $rest)
        } else {
            $out += $w
        }
    }
    return [string]::Join(" ", $out)
}

for ($r = 2; $r -le 1705; $r++) {
    foreach ($c in 1, 2) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null) {
            $newVal = TitleCaseWords $val
            # Use .Equals() (ordinal / case-sensitive) rather than -eq/-ne,
            # which this host treats as case-insensitive.
            if (-not $newVal.Equals($val)) {
                $cell.Value = $newVal
            }
        }
    }
}

# Drop the trailing footnote/metadata rows (1707:1711) so the sheet's used
# range shrinks back down to A1:D1705. (Row 1706 is already blank.)
$ws.Range("A1707:A1711").EntireRow.Delete()
